# Applies price/volume/coin updates per commit diff for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.377.72'
$ws.Range("E2").Value = '  +0.47%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.483.76'
$ws.Range("E3").Value = '  +0.60%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.61%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.02'
$ws.Range("E5").Value = '  +0.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.74'
$ws.Range("E6").Value = '  -1.13%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.18%  '

# Row 8
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.518.98'
$ws.Range("E9").Value = '  +1.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0975'
$ws.Range("E10").Value = '  -1.24%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.16'
$ws.Range("E12").Value = '  -2.08%  '

# Row 13
$ws.Range("E13").Value = '  -2.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.961.70'
$ws.Range("E14").Value = '  +1.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.605.94'
$ws.Range("E15").Value = '  +0.99%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.13'
$ws.Range("E16").Value = '  -1.54%  '

# Row 17
$ws.Range("E17").Value = '  -0.61%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.511.80'
$ws.Range("E18").Value = '  +1.73%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.66'
$ws.Range("E19").Value = '  -0.55%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '322.37'
$ws.Range("E20").Value = '  +0.21%  '

# Row 21
$ws.Range("E21").Value = '  -0.51%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.14'
$ws.Range("E22").Value = '  +6.72%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.45'
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("E25").Value = '  -1.67%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("E27").Value = '  +0.78%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.40'
$ws.Range("E28").Value = '  -0.19%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0754'
$ws.Range("E29").Value = '  +0.40%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  +1.46%  '

# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  +1.05%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.91'
$ws.Range("E32").Value = '  -0.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.26'
$ws.Range("E33").Value = '  -0.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.25%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.10'
$ws.Range("E36").Value = '  +0.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  -5.01%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("E38").Value = '  -1.00%  '

# Row 39
$ws.Range("E39").Value = '  -0.61%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.21'
$ws.Range("E40").Value = '  -0.38%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.778'
$ws.Range("E41").Value = '  -3.67%  '

# Row 42
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.50'
$ws.Range("E42").Value = '  +0.68%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '278.57'
$ws.Range("E43").Value = '  +1.72%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.02'
$ws.Range("E44").Value = '  -1.79%  '

# Row 45
$ws.Range("E45").Value = '  +1.31%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '122.89'
$ws.Range("E46").Value = '  -1.20%  '

# Row 47
$ws.Range("E47").Value = '  +0.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0501'
$ws.Range("E48").Value = '  +2.31%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.81'
$ws.Range("E49").Value = '  +0.53%  '

# Row 50
$ws.Range("E50").Value = '  +0.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.92'
$ws.Range("E51").Value = '  -0.74%  '
